$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.470.71'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.813.03'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.57'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.594'
$ws.Range('E6').Value = '  +2.95%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '38.30'
$ws.Range('E8').Value = '  +6.66%  '
$ws.Range('E9').Value = '  -3.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0675'
$ws.Range('E10').Value = '  -2.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0972'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.074.55'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.22'
$ws.Range('E13').Value = '  -2.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.814.38'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.468.96'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('E17').Value = '  -1.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.26'
$ws.Range('E18').Value = '  -1.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.11'
$ws.Range('E19').Value = '  -1.01%  '
$ws.Range('E20').Value = '  -2.67%  '
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('E24').Value = '  +3.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.19'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.81'
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.54'
$ws.Range('E27').Value = '  +4.00%  '
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.80'
$ws.Range('E30').Value = '  -1.14%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.23'
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('E32').Value = '  -2.58%  '
$ws.Range('E33').Value = '  -5.06%  '
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.361.30'
$ws.Range('E35').Value = '  -2.51%  '
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  -5.03%  '
$ws.Range('E39').Value = '  -1.64%  '
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '81.74'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.81'
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0508'
$ws.Range('E46').Value = '  +1.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.975.62'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.76'
$ws.Range('E48').Value = '  -4.54%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.04'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('E51').Value = '  -4.80%  '
